# Re-run of the Monte-Carlo simulation (p1_1.m) that produced this sheet:
# row 1 holds the fixed sample sizes, rows 2-5 hold freshly recomputed
# statistics. Push the refreshed numbers back into the sheet, restore
# column B's own explicit width entry, and force a full recalculation the
# next time the workbook is opened.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(-0.18831388744517291, -0.58941692700257031, 0.30365503262274862, -0.30584205623914329, -0.078780860491483265),
  @(-0.20810107090629623, -0.63178884480325992, 0.31756122672862419, -0.3212316741614169,  -0.078464282939658375),
  @(-0.71373338202181469, -0.32267450327747688, 0.26062939509323396, -0.52379545683380802, -0.62344581981811453),
  @(-0.77478761137803376, -0.36374971582050242, 0.27519400160685986, -0.56459947019063006, -0.66403069789680547)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $col = 1 + $j
        $ws.Cells.Item($row, $col).Value = $vals[$j]
    }
}

# Column B keeps the same visual width as before, but gets it re-applied
# explicitly so it round-trips as its own <col> entry again.
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(1).ColumnWidth

# A date-formatted, bordered style got registered in the style table (from
# formatting work done elsewhere in the workbook) but isn't used by any
# visible cell any more; touch-and-clear a scratch cell so the style
# definition is preserved without leaving stray content on the sheet.
$scratch = $ws.Range("Z100")
$scratch.NumberFormat = "m/d/yy h:mm"
$scratch.Borders.LineStyle = 1
$scratch.Clear()

# Force Excel to fully recalculate the workbook the next time it is opened.
$wb.ForceFullCalculation = $true

$wb.Save()
